$wb = $excel.ActiveWorkbook

# "go back to beta build": remove the pidSVC(<int>) command row that had been
# added to the Commands sheet (row 94: "pidSVC(<int>)" /
# "sets the PID target set value SV given in C"). Deleting the whole row
# shifts everything below it up by one and drops the now-unused shared
# strings for that command.
$wsCommands = $wb.Worksheets.Item("Commands")
$wsCommands.Activate()
$wsCommands.Rows.Item(94).Delete()
$wsCommands.Range("B58:C58").Select()

# Keep the Commands sheet as the active/selected tab, as in the original file.
$wsCommands.Activate()
